$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.930.26"
$ws.Range("E2").Value = "  -1.70%  "
$ws.Range("D3").Value = "'1.834.12"
$ws.Range("E3").Value = "  -1.87%  "
$ws.Range("D4").Value = "'0.9998"
$ws.Range("E4").Value = "  -0.24%  "
$ws.Range("D5").Value = "'245.25"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").Value = "'0.6906"
$ws.Range("E6").Value = "  -2.27%  "
$ws.Range("D7").Value = "'1.0000"
$ws.Range("E7").Value = "  -0.23%  "
$ws.Range("D8").Value = "'0.07671"
$ws.Range("E8").Value = "  -2.75%  "
$ws.Range("D9").Value = "'0.3051"
$ws.Range("E9").Value = "  -2.64%  "
$ws.Range("D10").Value = "'23.49"
$ws.Range("E10").Value = "  -4.36%  "
$ws.Range("D11").Value = "'0.07810"
$ws.Range("E11").Value = "  -0.97%  "
$ws.Range("D12").Value = "'1.833.95"
$ws.Range("E12").Value = "  -3.25%  "
$ws.Range("E13").Value = "  -2.43%  "
$ws.Range("D14").Value = "'90.40"
$ws.Range("E14").Value = "  -3.26%  "
$ws.Range("D15").Value = "'0.6796"
$ws.Range("E15").Value = "  -2.98%  "
$ws.Range("D16").Value = "'6.423"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").Value = "'0.000008335"
$ws.Range("E17").Value = "  -0.18%  "
$ws.Range("D18").Value = "'28.928.31"
$ws.Range("E18").Value = "  -1.95%  "
$ws.Range("D19").Value = "'242.86"
$ws.Range("E19").Value = "  -3.58%  "
$ws.Range("D20").Value = "'2.082.11"
$ws.Range("E20").Value = "  -3.08%  "
$ws.Range("E21").Value = "  -3.21%  "
$ws.Range("D22").Value = "'1.0000"
$ws.Range("E22").Value = "  -0.19%  "
$ws.Range("D23").Value = "'7.457"
$ws.Range("E23").Value = "  -2.24%  "
$ws.Range("D24").Value = "'1.000"
$ws.Range("E24").Value = "  -0.22%  "
$ws.Range("E25").Value = "  -5.76%  "
$ws.Range("D26").Value = "'161.02"
$ws.Range("E26").Value = "  -0.61%  "
$ws.Range("D27").Value = "'8.794"
$ws.Range("E27").Value = "  -2.23%  "
$ws.Range("E28").Value = "  -2.88%  "
$ws.Range("D29").Value = "'1.560"
$ws.Range("E29").Value = "  +3.42%  "
$ws.Range("D30").Value = "'4.213"
$ws.Range("D31").Value = "'4.147"
$ws.Range("E31").Value = "  -2.52%  "
$ws.Range("D32").Value = "'1.178"
$ws.Range("E32").Value = "  -2.40%  "
$ws.Range("D33").Value = "'0.05111"
$ws.Range("E33").Value = "  -3.78%  "
$ws.Range("D34").Value = "'0.7581"
$ws.Range("E34").Value = "  +1.10%  "
$ws.Range("D35").Value = "'1.846"
$ws.Range("E35").Value = "  -2.55%  "
$ws.Range("D36").Value = "'1.147"
$ws.Range("E36").Value = "  -2.08%  "
$ws.Range("D37").Value = "'2.681"
$ws.Range("E37").Value = "  -1.33%  "
$ws.Range("D38").Value = "'0.01842"
$ws.Range("E38").Value = "  -2.36%  "
$ws.Range("D39").Value = "'1.232.59"
$ws.Range("E39").Value = "  -3.44%  "
$ws.Range("D40").Value = "'2.694"
$ws.Range("E40").Value = "  -2.64%  "
$ws.Range("D41").Value = "'0.9230"
$ws.Range("E41").Value = "  +3.42%  "
$ws.Range("D42").Value = "'108.74"
$ws.Range("E42").Value = "  -0.13%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.774"
$ws.Range("E43").Value = "  -4.96%  "
$ws.Range("B44").Value = "PaxDollar"
$ws.Range("C44").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D44").Value = "'0.9995"
$ws.Range("E44").Value = "  -0.24%  "
$ws.Range("B45").Value = "Mantle"
$ws.Range("C45").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D45").Value = "'0.5174"
$ws.Range("E45").Value = "  -0.28%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "'1.981.57"
$ws.Range("E46").Value = "  -3.35%  "
$ws.Range("B47").Value = "BabyDogeCoin"
$ws.Range("C47").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D47").Value = "'0.00000000122"
$ws.Range("E47").Value = "  -4.69%  "
$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "'9.509"
$ws.Range("E48").Value = "  -0.47%  "
$ws.Range("D49").Value = "'63.89"
$ws.Range("E49").Value = "  -10.34%  "
$ws.Range("E50").Value = "  -3.02%  "
$ws.Range("D51").Value = "'6.910"
